$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "January " + "7" (two runs) -> "January 7" (single run), mirroring
#    the pattern already used by the other "January Nth:" headings.
# ------------------------------------------------------------------
$headingRange = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("January 7")) {
        $headingRange = $para.Range
        break
    }
}

$start = $headingRange.Start
$mergeRange = $d.Range($start, $start + 9)   # covers "January 7"
$mergeRange.Delete()
$insertion = $d.Range($start, $start)
$insertion.InsertBefore("January 7")

# ------------------------------------------------------------------
# 2) After the paragraph ending in the loader-cutting-pictures ";",
#    add a blank paragraph plus the new "January 9th:" entry and its
#    write-up.
# ------------------------------------------------------------------
$anchorPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*cuts up a picture into rectangles of a desired size;*") {
        $anchorPara = $para
        break
    }
}

$insertPos = $anchorPara.Range.End - 1   # just before this paragraph's own mark
$insertRange = $d.Range($insertPos, $insertPos)

$newParasXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>' + `
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:r><w:lastRenderedPageBreak/><w:t>January 9</w:t></w:r>' + `
        '<w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>th</w:t></w:r>' + `
        '<w:r><w:t>:</w:t></w:r>' + `
    '</w:p>' + `
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:r><w:t>Fixed issue with loading an asset</w:t></w:r>' + `
        '<w:r><w:t>, started work on system where I can put pngs in a folder and depending on the name of the folder the assets will be loaded in and give a certain attribure. Right have all the assets in a folder loading in. need to make manager which will create a tile depending with relevant information</w:t></w:r>' + `
        '<w:r><w:t>.</w:t></w:r>' + `
    '</w:p>'

[void]$insertRange.InsertXML($newParasXml)
